$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (row 1 title)
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 14:16"

# Full country table (re-pulled data + re-sorted by "Casos totales" desc),
# encoded as pipe-delimited rows: Pais|Casos totales|Nuevos casos|Casos activos|Recuperados|Casos criticos|Muertes hoy|Muertes
$data = @"
China|81054|46|72440|5353|1845|6|3261
Italia|53578|0|6072|42681|2857|0|4825
España|28572|3076|2125|24722|1785|344|1725
Estados Unidos|26906|2699|178|26380|708|46|348
Alemania|23129|765|209|22827|2|9|93
Iran|21638|1028|7635|12318|0|129|1685
Francia|14459|0|1587|12310|1525|0|562
Corea del Sur|8897|98|2909|5884|59|2|104
Suiza|7225|362|131|7014|141|0|80
Reino Unido|5018|0|93|4685|20|7|240
Paises Bajos|3631|0|2|3493|354|0|136
Belgica|3401|586|263|3063|288|8|75
Austria|3214|222|9|3196|15|1|9
Noruega|2246|82|6|2233|28|0|7
Suecia|1906|136|16|1869|68|1|21
Portugal|1600|320|5|1581|26|2|14
Dinamarca|1395|69|1|1381|42|0|13
Australia|1353|281|46|1300|2|0|7
Canada|1331|3|14|1297|1|1|20
Malasia|1306|123|139|1157|26|2|10
Brasil|1201|23|2|1181|18|0|18
Japon|1054|0|215|803|55|0|36
Chequia|1047|52|6|1041|7|0|0
Turquia|947|0|0|926|0|0|21
Israel|945|62|37|907|20|0|1
Luxemburgo|798|128|6|784|3|0|8
Irlanda|785|0|5|777|13|0|3
Pakistan|730|85|13|714|0|0|3
Crucero|712|0|567|137|15|0|8
Finlandia|626|103|10|615|12|0|1
Tailandia|599|188|44|554|7|0|1
Islandia|568|95|5|562|1|0|1
Polonia|563|27|13|543|3|2|7
Chile|537|0|8|528|7|0|1
Ecuador|532|0|3|522|2|0|7
Grecia|530|0|19|498|18|0|13
Indonesia|514|64|29|437|0|10|48
Arabia Saudita|511|119|17|494|0|0|0
Catar|481|0|27|454|6|0|0
Singapur|455|23|144|309|14|0|2
Rumania|433|66|64|367|14|2|2
Eslovenia|383|0|0|382|12|0|1
Filipinas|380|73|15|340|1|6|25
India|360|28|24|329|0|2|7
Barein|332|22|149|182|4|0|1
Estonia|326|20|4|322|2|0|0
Peru|318|0|1|312|5|0|5
Rusia|306|0|16|289|0|0|1
Egipto|294|0|42|242|0|0|10
Hong Kong|274|0|100|170|4|0|4
Mexico|251|48|4|245|1|0|2
Libano|248|18|8|236|4|0|4
Panama|245|0|1|241|7|0|3
Sudafrica|240|0|2|238|0|0|0
Croacia|235|29|5|229|5|0|1
Irak|233|19|57|156|0|3|20
Argentina|225|67|27|194|0|0|4
Colombia|210|14|3|206|0|1|1
Armenia|190|30|2|188|6|0|0
Serbia|188|17|2|184|4|1|2
Kuwait|188|12|27|161|5|0|0
Eslovaquia|178|0|7|171|2|0|0
San Marino|175|15|4|151|13|0|20
Bulgaria|174|11|3|168|3|0|3
Taiwan|169|16|28|139|0|0|2
Emiratos Arabes Unidos|153|0|38|113|2|0|2
Letonia|139|15|1|138|0|0|0
Argelia|139|0|65|59|0|0|15
Uruguay|135|25|0|135|2|0|0
Hungria|131|28|16|109|6|2|6
Costa Rica|117|0|2|113|2|0|2
Republica de Macedonia|114|29|1|113|1|0|0
Republica Dominicana|112|0|0|109|0|0|3
Marruecos|108|12|3|102|1|0|3
Lituania|105|6|1|103|1|0|1
Jordania|100|0|1|99|0|0|0
Vietnam|99|5|17|82|2|0|0
Bosnia y Herzegovina|93|0|2|90|1|0|1
Islas Feroe|92|0|3|89|0|0|0
Albania|89|13|2|85|2|0|2
Principado de Andorra|88|0|1|86|2|1|1
Republica de Chipre|84|0|3|80|3|0|1
Brunei|83|0|1|82|2|0|0
Sri Lanka|82|5|3|79|2|0|0
Moldavia|80|0|1|78|3|0|1
Bielorrusia|76|0|15|61|0|0|0
Tunez|75|15|1|71|7|2|3
Burkina Faso|75|11|5|66|0|1|4
Malta|73|0|2|71|1|0|0
Venezuela|70|0|15|55|2|0|0
Nueva Zelanda|66|14|0|66|0|0|0
Estado de Palestina|59|6|17|42|0|0|0
Kazajistan|57|3|0|57|0|0|0
Guadalupe|56|0|0|55|4|0|1
Senegal|56|0|5|51|0|0|0
Oman|55|3|13|42|0|0|0
Georgia|54|5|1|53|1|0|0
Camboya|53|0|2|51|0|0|0
Azerbaiyan|53|0|11|41|0|0|1
Trinidad yTobago|49|0|0|49|0|0|0
Reunion|47|0|0|47|0|0|0
Ucrania|47|0|1|43|0|0|3
Uzbekistan|42|1|0|42|0|0|0
Camerun|40|0|2|38|0|0|0
Liechtenstein|37|0|0|37|0|0|0
Martinica|37|0|0|36|7|0|1
Afganistan|34|10|1|33|0|0|0
Guam|27|12|0|26|0|1|1
Banglades|27|3|3|22|0|0|2
Honduras|26|2|0|26|0|0|0
Nigeria|26|4|1|25|0|0|0
Puerto Rico|23|2|0|22|0|0|1
Consejo Danes para los Refugiados|23|0|0|22|0|0|1
Paraguay|22|0|0|21|1|0|1
Cuba|21|0|0|20|0|0|1
Ghana|21|0|0|20|0|0|1
Macao|20|1|10|10|0|0|0
Bolivia|19|0|0|19|0|0|0
Jamaica|19|0|2|16|0|0|1
Guayana Francesa|18|0|0|18|0|0|0
Mauricio|18|4|0|17|0|0|1
Monaco|18|0|1|17|0|0|0
Guyana|18|0|0|17|0|0|1
Ruanda|17|0|0|17|0|0|0
Guatemala|17|0|0|16|0|0|1
Togo|16|0|0|16|0|0|0
Montenegro|16|0|0|16|0|0|0
Polinesia Francesa|15|0|0|15|0|0|0
Kirguistan|14|0|0|14|0|0|0
Barbados|14|0|0|14|0|0|0
Costa de Marfil|14|0|1|13|0|0|0
Maldivas|13|0|3|10|0|0|0
Tanzania|12|6|0|12|0|0|0
Etiopia|11|2|0|11|0|0|0
Mayotte|11|0|0|11|0|0|0
Mongolia|10|0|0|10|0|0|0
Gibraltar|10|0|2|8|0|0|0
Aruba|8|3|1|7|0|0|0
Kenia|7|0|0|7|0|0|0
Seychelles|7|0|0|7|0|0|0
Islas Virgenes de los Estados Unidos|6|0|0|6|0|0|0
Guinea Ecuatorial|6|0|0|6|0|0|0
Surinam|5|0|0|5|0|0|0
San Martin (Parte Francesa)|5|0|0|5|0|0|0
Isla de Man|5|3|0|5|0|0|0
Gabon|5|0|0|4|0|0|1
Nueva Caledonia|4|0|0|4|0|0|0
Suazilandia|4|3|0|4|0|0|0
Bahamas|4|0|0|4|0|0|0
Liberia|3|0|0|3|0|0|0
Madagascar|3|0|0|3|0|0|0
Zambia|3|1|0|3|0|0|0
San Bartolome|3|0|0|3|0|0|0
Zimbabue|3|0|0|3|0|0|0
Republica de Africa Central|3|0|0|3|0|0|0
Namibia|3|0|0|3|0|0|0
Cabo Verde|3|0|0|3|0|0|0
El Salvador|3|0|0|3|0|0|0
Congo|3|0|0|3|0|0|0
Islas Caimanes|3|0|0|2|0|0|1
Curazao|3|0|0|2|0|0|1
Haiti|2|0|0|2|0|0|0
Butan|2|0|0|2|0|0|0
Groenlandia|2|0|0|2|0|0|0
Fiyi|2|0|0|2|0|0|0
Santa Lucia|2|0|0|2|0|0|0
Guinea|2|0|0|2|0|0|0
Benin|2|0|0|2|0|0|0
Bermudas|2|0|0|2|0|0|0
Nicaragua|2|0|0|2|0|0|0
Mauritania|2|0|0|2|0|0|0
Angola|2|0|0|2|0|0|0
Sudan|2|0|0|1|0|0|1
San Martin (Parte Holandesa)|1|0|0|1|0|0|0
Montserrat|1|0|0|1|0|0|0
Gambia|1|0|0|1|0|0|0
San Vicente y las Granadinas|1|0|0|1|0|0|0
Santa Sede|1|0|0|1|0|0|0
Republica del Chad|1|0|0|1|0|0|0
Niger|1|0|0|1|0|0|0
Eritrea|1|0|0|1|0|0|0
Papua Nueva Guinea|1|0|0|1|0|0|0
Uganda|1|0|0|1|0|0|0
Somalia|1|0|0|1|0|0|0
Republica de Yibuti|1|0|0|1|0|0|0
Timor Oriental|1|0|0|1|0|0|0
Antigua y Barbuda|1|0|0|1|0|0|0
Nepal|1|0|1|0|0|0|0
"@

$lines = $data -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Length
$colCount = 8

$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $lines[$i] -split '\|'
    $arr[$i,0] = $parts[0]
    for ($j = 1; $j -lt $colCount; $j++) {
        $arr[$i,$j] = [double]$parts[$j]
    }
}

$startRow = 4
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range("A" + $startRow + ":H" + $endRow)
$rng.Value2 = $arr
